$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "When valid bill details are provided 0400005777052"
$ws.Range("F2").Value = "0400005777052"

$newSelectQuery = "SELECT LB.BILL_AMOUNT, LB.COMPANY_CODE, LB.DUE_DATE FROM LP_BILLS LB WHERE LB.CONSUMER_NO='{ConsumerNo}' ORDER BY LB.CREATED_ON DESC"
$ws.Range("G2").Value = $newSelectQuery
$ws.Range("G3").Value = $newSelectQuery
$ws.Range("G4").Value = $newSelectQuery
$ws.Range("G5").Value = $newSelectQuery

$newBeginUpdate = "BEGIN UPDATE DC_SCHEDULED_TRAN_MASTER STM SET STM.STATE = 46 , STM.IS_DELETED = 1 WHERE STM.BILL_BENEFICIARY_ID = (SELECT BPB.BENEFICIARY_ID FROM DC_BILL_PAYMENT_BENEFICIARY BPB WHERE BPB.CONSUMER_NUMBER = '{ConsumerNo}' AND BPB.CUSTOMER_INFO_ID = (SELECT CI.CUSTOMER_INFO_ID FROM DC_CUSTOMER_INFO CI WHERE CI.CUSTOMER_NAME = '{customer_name}') AND BPB.IS_ACTIVE = 1);UPDATE DC_BILL_PAYMENT_BENEFICIARY DPB SET DPB.IS_SI_SCHEDULED = 0,DPB.IS_ACTIVE = 0 WHERE DPB.CONSUMER_NUMBER = '{ConsumerNo}' AND DPB.CUSTOMER_INFO_ID = (SELECT CI.CUSTOMER_INFO_ID FROM DC_CUSTOMER_INFO CI WHERE CI.CUSTOMER_NAME = '{customer_name}') AND DPB.IS_ACTIVE = 1;COMMIT;END;"
$ws.Range("C2").Value = $newBeginUpdate
$ws.Range("C3").Value = $newBeginUpdate
$ws.Range("C4").Value = $newBeginUpdate
$ws.Range("C5").Value = $newBeginUpdate

$ws.Columns.Item(7).ColumnWidth = 140.42578125

$ws.Range("D13").Select()
